$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A15").Value = "Oakwood Village Toronto"
$ws.Range("B15").Value = 43.7491
$ws.Range("C15").Value = -79.4404
$ws.Range("D15").Value = 1.0406060606060596

$ws.Range("A16").Value = "Vandenyoung City"
$ws.Range("B16").Value = 43.6757
$ws.Range("C16").Value = -79.4181
$ws.Range("D16").Value = 6.5175595238095223

$ws.Range("A17").Value = "Moore Park"
$ws.Range("B17").Value = 43.7519
$ws.Range("C17").Value = -79.3834
$ws.Range("D17").Value = 5.6800595238095228
